# "10Th - MB for single stock and added new group"
#
# This MarketBeat-rank watch-sheet tracks, for each analyst/firm (rows),
# whether their rating changed on a given date (columns, newest date
# left-most). This edit:
#   1. Adds two brand-new dates: Jun_27 and Jun_26 (Jun_26 gets TWO
#      columns because two different upgrade notes landed that day).
#   2. Records a "Hold -> Buy" upgrade for Jefferies Financial Group
#      (row 13) on 6/26/2018, highlighted with a fill color.
#   3. Adds a new group: two more firms/benchmarks ("Benchmark" and
#      "Evercore ISI") as new rows at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room for the two new dates -----------------------------
# Insert three fresh columns before column B (the old "newest date"
# column). This pushes the existing B:E ("Jun_17","Jun_15","Jun_13",
# "Jun_10") right to E:H, exactly like Excel's own Insert does.
$ws.Columns("B:D").Insert()

# Re-apply a consistent custom width across the whole date block
# (C:H), matching the look of the rest of the sheet.
$ws.Columns("C:H").ColumnWidth = 8.0

# --- 2. New date headers ---------------------------------------------
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# --- 3. Default every analyst row to "UN" (unchanged) for the new ----
#        Jun_27 / Jun_26 / Jun_26 columns.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# --- 4. Record the single-stock upgrade for Jefferies Financial ------
#        Group (row 13), on 6/26/2018, and highlight it.
$upgradeNote = "6/26/2018,Upgrades,Hold -> Buy,"
$ws.Range("B13").Value = $upgradeNote
$ws.Range("C13").Value = $upgradeNote
$ws.Range("D13").Value = $upgradeNote
$ws.Range("B13:D13").Interior.ColorIndex = 42

# --- 5. Added new group: two more rows at the bottom -----------------
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
